$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, [string]$text)
    # Force the cell to be treated as text so numeric-looking
    # strings (e.g. "86.10", "2.00") keep their exact formatting
    # instead of being coerced into numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
$ws.Cells.Item(2, 4).Value = '26.883.69'
$ws.Cells.Item(2, 5).Value = '  -1.38%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.563.97'
$ws.Cells.Item(3, 5).Value = '  +0.11%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.39%  '

# Row 5
Set-TextCell $ws.Cells.Item(5, 4) '205.98'
$ws.Cells.Item(5, 5).Value = '  -0.01%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  -1.26%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.32%  '

# Row 8
Set-TextCell $ws.Cells.Item(8, 4) '21.78'
$ws.Cells.Item(8, 5).Value = '  -1.19%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -0.19%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -1.11%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +0.30%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.785.82'
$ws.Cells.Item(12, 5).Value = '  -0.05%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '1.561.63'
$ws.Cells.Item(13, 5).Value = '  -0.16%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  -1.13%  '

# Row 15
Set-TextCell $ws.Cells.Item(15, 4) '0.514'
$ws.Cells.Item(15, 5).Value = '  -0.23%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '26.855.35'
$ws.Cells.Item(16, 5).Value = '  -1.36%  '

# Row 17
Set-TextCell $ws.Cells.Item(17, 4) '61.16'
$ws.Cells.Item(17, 5).Value = '  -3.12%  '

# Row 18
$ws.Cells.Item(18, 2).Value = 'Chainlink'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell $ws.Cells.Item(18, 4) '7.36'
$ws.Cells.Item(18, 5).Value = '  +2.20%  '

# Row 19
$ws.Cells.Item(19, 2).Value = 'BitcoinCash'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell $ws.Cells.Item(19, 4) '214.25'
$ws.Cells.Item(19, 5).Value = '  +1.62%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '0.0₃0679'
$ws.Cells.Item(20, 5).Value = '  -1.22%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  -0.37%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +0.61%  '

# Row 23
Set-TextCell $ws.Cells.Item(23, 4) '9.19'
$ws.Cells.Item(23, 5).Value = '  -2.42%  '

# Row 24
Set-TextCell $ws.Cells.Item(24, 4) '2.00'
$ws.Cells.Item(24, 5).Value = '  -0.07%  '

# Row 25
Set-TextCell $ws.Cells.Item(25, 4) '154.01'
$ws.Cells.Item(25, 5).Value = '  +0.82%  '

# Row 26
Set-TextCell $ws.Cells.Item(26, 4) '6.71'
$ws.Cells.Item(26, 5).Value = '  +1.39%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +0.58%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -0.45%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -1.01%  '

# Row 30
Set-TextCell $ws.Cells.Item(30, 4) '0.0463'
$ws.Cells.Item(30, 5).Value = '  -0.66%  '

# Row 31
Set-TextCell $ws.Cells.Item(31, 4) '1.11'
$ws.Cells.Item(31, 5).Value = '  -3.26%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +0.07%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '1.403.15'
$ws.Cells.Item(33, 5).Value = '  +1.61%  '

# Row 34
Set-TextCell $ws.Cells.Item(34, 4) '2.92'
$ws.Cells.Item(34, 5).Value = '  -0.84%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -1.09%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -1.45%  '

# Row 37
Set-TextCell $ws.Cells.Item(37, 4) '0.921'
$ws.Cells.Item(37, 5).Value = '  -1.99%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +0.04%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +0.25%  '

# Row 40
Set-TextCell $ws.Cells.Item(40, 4) '0.814'
$ws.Cells.Item(40, 5).Value = '  +0.18%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.39%  '

# Row 42
Set-TextCell $ws.Cells.Item(42, 4) '0.996'
$ws.Cells.Item(42, 5).Value = '  +0.18%  '

# Row 43
Set-TextCell $ws.Cells.Item(43, 4) '5.33'
$ws.Cells.Item(43, 5).Value = '  +2.15%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'RenderToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws.Cells.Item(44, 4) '1.77'
$ws.Cells.Item(44, 5).Value = '  -1.86%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'MXToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell $ws.Cells.Item(45, 4) '2.17'
$ws.Cells.Item(45, 5).Value = '  -3.88%  '

# Row 46
Set-TextCell $ws.Cells.Item(46, 4) '63.14'
$ws.Cells.Item(46, 5).Value = '  -0.45%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '1.698.90'
$ws.Cells.Item(47, 5).Value = '  -0.25%  '

# Row 48
Set-TextCell $ws.Cells.Item(48, 4) '86.10'
$ws.Cells.Item(48, 5).Value = '  +0.93%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Cronos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws.Cells.Item(49, 4) '0.0505'
$ws.Cells.Item(49, 5).Value = '  +2.39%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(50, 4).Value = '0.0₇0985'
$ws.Cells.Item(50, 5).Value = '  -1.17%  '

# Row 51
Set-TextCell $ws.Cells.Item(51, 4) '0.0947'
$ws.Cells.Item(51, 5).Value = '  +0.58%  '

Write-Host "Applied cryptos.xlsx updates"